$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# Row 2 formulas
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("H2").Formula = "=SUM(G2:G11)"
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# Row 3 formula (own, not part of the shared group)
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

# Rows 4-15 share the same relative formula pattern (D_row - D_row-1) * B_row / 100
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# Match final selection left by the author
$ws.Range("A2").Select() | Out-Null

$wb.Save()
